$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in "Wie" (D) / "Opmerkingen" (E) text values for various rows ---
$ws.Range("D4").Value = "Laura"
$ws.Range("D5").Value = "Laura"

$ws.Range("D8").Value = "Nicole"
$ws.Range("E8").Value = "in geschiedenis"

$ws.Range("E9").Value = "?"

$ws.Range("D10").Value = "Nicole"
$ws.Range("D11").Value = "Laura"
$ws.Range("D12").Value = "Nicole"

$ws.Range("D14").Value = "Laura"

$ws.Range("D16").Value = "nicole/Laura"

$ws.Range("D18").Value = "nicole/Laura"

$ws.Range("D20").Value = "Laura/nicole"

$ws.Range("D23").Value = "Laura"
$ws.Range("D24").Value = "Nicole"
$ws.Range("D25").Value = "Laura"

$ws.Range("D27").Value = "Laura"

# --- Add new empty yellow-highlighted cells (same style as the existing B3/C19 markers) ---
$ws.Range("B13").Interior.Color = 65535
$ws.Range("C17").Interior.Color = 65535
$ws.Range("C18").Interior.Color = 65535
$ws.Range("C20").Interior.Color = 65535
$ws.Range("C24").Interior.Color = 65535

# --- The old yellow marker in B24 moved to C24, so clear B24 entirely ---
$ws.Range("B24").Clear()

# --- Update the view: scroll back to top-left and select D1 instead of B24 ---
[void]$ws.Range("D1").Select()
